$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2;  A=94177; B="Aylla Cardoso";        C="TI";                  D="Viagem de negocios"; E=7; F=45100; G=4198.8 },
    @{ Row=3;  A=2214;  B="Maitê da Luz";          C="Recursos Humanos";    D="Consulta medica";    E=6; F=45086; G=2055.33 },
    @{ Row=4;  A=55363; B="Srta. Liz Pereira";     C="TI";                  D="Viagem de negocios"; E=2; F=45103; G=7773.59 },
    @{ Row=5;  A=28003; B="João Vitor Fernandes";  C="P&D";                 D="Outros";             E=6; F=45101; G=2873.4 },
    @{ Row=6;  A=45721; B="Manuela Castro";        C="Recursos Humanos";    D="Problemas pessoais"; E=8; F=45098; G=6390.97 },
    @{ Row=7;  A=44110; B="Alice Câmara";          C="Recursos Humanos";    D="Problemas pessoais"; E=2; F=45104; G=2386.39 },
    @{ Row=8;  A=41952; B="Antônio Ribeiro";       C="Financeiro";          D="Viagem de negocios"; E=2; F=45082; G=9196.459999999999 },
    @{ Row=9;  A=62586; B="Kevin Vieira";          C="P&D";                 D="Problemas pessoais"; E=6; F=45106; G=3231.66 },
    @{ Row=10; A=12970; B="Luara Sá";              C="Engenharia";          D="Consulta medica";    E=3; F=45100; G=4866.08 },
    @{ Row=11; A=65106; B="Henrique da Rosa";      C="Marketing";           D="Problemas pessoais"; E=4; F=45085; G=4426.97 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.A
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
}
